# Contact Import Updates (2023-08-12)
#
# - Drop the unused "Sheet2" / "Sheet3" tabs, leaving only the contacts
#   header sheet ("Sheet1").
# - Remove the "Reports To" and "Assigned To" columns from the CSV header
#   row (columns T:U), shifting every column to their right left by two.
# - Refresh the saved selection/active cell to reflect where the user was
#   working when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Reports To" (T1) and "Assigned To" (U1) columns; this shifts
# the remaining headers (Email OPT OUT onward) two columns to the left and
# shrinks the used range from A1:AI1 down to A1:AG1.
$ws.Columns("T:U").Delete()

# The workbook only ships the contacts header sheet going forward.
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Leave the selection where it was when the file was saved.
$ws.Range("V10").Select()
